$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text. Columns D/E in this sheet store
# plain-text numbers/percentages (e.g. "0.999", "12.30", "  +3.71%  ");
# a bare .Value assignment would let Excel auto-coerce numeric-looking
# strings into real numbers (dropping formatting like trailing zeros), so
# we briefly force the cell to Text format, write the string, then put the
# original (General) number format back.
function Set-TextValue($cell, $value) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = $fmt
}

$ws.Range("D2").Value = "56.376.96"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "2.972.67"
$ws.Range("E3").Value = "  +2.91%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "502.02"
$ws.Range("E5").Value = "  +5.12%  "
Set-TextValue $ws.Range("D6") "134.61"
$ws.Range("E6").Value = "  +6.22%  "
$ws.Range("E7").Value = "  -0.09%  "
Set-TextValue $ws.Range("D8") "0.428"
$ws.Range("E8").Value = "  +5.66%  "
Set-TextValue $ws.Range("D9") "7.41"
$ws.Range("E9").Value = "  +10.86%  "
Set-TextValue $ws.Range("D10") "0.107"
$ws.Range("E10").Value = "  +9.30%  "
$ws.Range("E11").Value = "  +4.59%  "
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("D13").Value = "3.475.00"
$ws.Range("E13").Value = "  +2.77%  "
Set-TextValue $ws.Range("D14") "25.26"
$ws.Range("E14").Value = "  +10.53%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D15") "0.0000151"
$ws.Range("E15").Value = "  +10.88%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "56.320.29"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("D17").Value = "2.966.12"
$ws.Range("E17").Value = "  +2.49%  "
Set-TextValue $ws.Range("D18") "5.71"
$ws.Range("E18").Value = "  +8.87%  "
Set-TextValue $ws.Range("D19") "12.30"
$ws.Range("E19").Value = "  +5.69%  "
Set-TextValue $ws.Range("D20") "7.71"
$ws.Range("E20").Value = "  +8.32%  "
Set-TextValue $ws.Range("D21") "321.82"
$ws.Range("E21").Value = "  +3.61%  "
$ws.Range("E22").Value = "  +0.24%  "
Set-TextValue $ws.Range("D23") "0.467"
$ws.Range("E23").Value = "  +3.92%  "
Set-TextValue $ws.Range("D24") "61.67"
$ws.Range("E24").Value = "  +3.04%  "
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  +5.20%  "
$ws.Range("D27").Value = "0.0₃0884"
$ws.Range("E27").Value = "  +7.08%  "
Set-TextValue $ws.Range("D28") "6.39"
$ws.Range("E28").Value = "  +1.57%  "
Set-TextValue $ws.Range("D29") "6.73"
$ws.Range("E29").Value = "  +7.80%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D30") "1.74"
$ws.Range("E30").Value = "  +7.15%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D31") "1.17"
$ws.Range("E31").Value = "  +2.09%  "
Set-TextValue $ws.Range("D32") "20.29"
$ws.Range("E32").Value = "  +5.64%  "
Set-TextValue $ws.Range("D33") "158.38"
$ws.Range("E33").Value = "  +13.04%  "
Set-TextValue $ws.Range("D34") "4.42"
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("E35").Value = "  +2.81%  "
Set-TextValue $ws.Range("D36") "5.52"
$ws.Range("E36").Value = "  +0.52%  "
Set-TextValue $ws.Range("D37") "0.0669"
$ws.Range("E37").Value = "  +7.40%  "
Set-TextValue $ws.Range("D38") "22.83"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "2.999.85"
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("E40").Value = "  -0.11%  "
Set-TextValue $ws.Range("D41") "36.10"
$ws.Range("E41").Value = "  +1.88%  "
Set-TextValue $ws.Range("D42") "0.637"
$ws.Range("E42").Value = "  +5.74%  "
$ws.Range("D43").Value = "2.231.16"
$ws.Range("E43").Value = "  +8.11%  "
Set-TextValue $ws.Range("D44") "1.39"
$ws.Range("E44").Value = "  +4.64%  "
Set-TextValue $ws.Range("D45") "0.976"
$ws.Range("E45").Value = "  +1.08%  "
Set-TextValue $ws.Range("D46") "3.54"
$ws.Range("E46").Value = "  +2.96%  "
Set-TextValue $ws.Range("D47") "1.92"
$ws.Range("E47").Value = "  +18.30%  "
$ws.Range("E48").Value = "  +9.71%  "
$ws.Range("E49").Value = "  +6.86%  "
Set-TextValue $ws.Range("D50") "18.86"
$ws.Range("E50").Value = "  +4.27%  "
$ws.Range("E51").Value = "  +8.03%  "
